$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: new "Fly back diode" line ---
# B12: price (1.21) with the same currency format as the rest of column B
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").Value = 1.21

# C12: quantity
$ws.Range("C12").Value = 15

# E12: hyperlink cell - give it the Hyperlink look (copy formats first so the
# eventual style lines up with the existing hyperlink cells), then register
# the actual hyperlink (this also writes the display text into the cell).
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E12"), "https://www.mouser.be/ProductDetail/ROHM-Semiconductor/1SS400CMT2R?qs=sGAEpiMZZMtoHjESLttvkn%252BvjfD1a1Smq9%2FW6eNwDXWDigaLofBvqg%3D%3D") | Out-Null
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null

# --- Column E width grows to fit the new (long) URL ---
$ws.Columns("E").ColumnWidth = 151.140625

# --- Selection moves to E15 ---
$ws.Range("E15").Select() | Out-Null

$excel.CutCopyMode = $false
